{"js": "// Append \"Coucou\" as a new run at the end of the document body's\n// (sole, currently empty) paragraph \u2014 right after the existing\n// bookmarkStart/bookmarkEnd pair.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Coucou\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append \"Coucou\" as a new run at the end of the document body's\n# (sole, currently empty) paragraph \u2014 right after the existing\n# bookmarkStart/bookmarkEnd pair.\n$d = $word.ActiveDocument\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$lastParagraph.Range.InsertAfter(\"Coucou\")\n"}
